# Plot of state/transition map
# - Added code to make a plot of a group of states and all the transitions
#   involving those states.

$wb = $excel.ActiveWorkbook

# --- General sheet: correct target-personnel / sim-length scale (25000 -> 2500, 100 -> 50) ---
$wsGeneral = $wb.Worksheets.Item("General")
$wsGeneral.Range("B5").Value = 2500
$wsGeneral.Range("B7").Value = 50

# --- Recruitment sheet: correct min/max recruitment age scale (/10) ---
$wsRecruitment = $wb.Worksheets.Item("Recruitment")
$wsRecruitment.Range("F8").Value = 35
$wsRecruitment.Range("F9").Value = 55
$wsRecruitment.Range("B17").Value = 35
$wsRecruitment.Range("B18").Value = 45
$wsRecruitment.Range("B19").Value = 55

# --- Output plots sheet: one fewer subpopulation graph (3 -> 2) now that the ---
# --- state/transition map plot covers that case separately ---
$wsOutput = $wb.Worksheets.Item("Output plots")
$wsOutput.Range("F3").Value = 2

# --- Make "Output plots" the active sheet/tab with F4 selected ---
$wsOutput.Activate() | Out-Null
$wsOutput.Range("F4").Select() | Out-Null
